$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.750.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.293.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +17.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.04%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.92'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.29%  '
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.638.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.859'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.289.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.703.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("E20").Value = '  +11.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("E22").Value = '  -1.76%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +10.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.97%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.38%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '175.44'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0925'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.13%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.71'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0362'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.28%  '
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.92'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +14.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.84%  '
$ws.Range("B43").Value = 'LidoDAOToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.241'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +22.53%  '
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("E47").Value = '  +3.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.38'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.72%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0994'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("E51").Value = '  +1.69%  '
